{"js": "// Remove the trailing \"Ver no Jupiter ...\" and \"\u00a9 2020 ...\" footer\n// paragraphs from the document body, leaving the surrounding empty\n// paragraphs (and everything else) untouched.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\nconst targetTexts = [\n  \"Ver no Jupiter Salvar em pdf Salvar em docx\",\n  \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n];\n\n// Walk backwards so deleting items doesn't disturb indices we still need.\nfor (let i = paragraphs.items.length - 1; i >= 0; i--) {\n  const paragraph = paragraphs.items[i];\n  if (targetTexts.indexOf(paragraph.text) !== -1) {\n    paragraph.delete();\n  }\n}\n\nawait context.sync();\n", "ps1": "# Remove the trailing \"Ver no Jupiter ...\" and \"\u00a9 2020 ...\" footer\n# paragraphs from the document, leaving the surrounding empty\n# paragraphs (and everything else) untouched.\n$d = $word.ActiveDocument\n\n$targets = @(\n  \"Ver no Jupiter Salvar em pdf Salvar em docx\",\n  \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n)\n\n$toDelete = New-Object System.Collections.ArrayList\nforeach ($p in $d.Paragraphs) {\n  $text = $p.Range.Text.TrimEnd([char]13, [char]7)\n  if ($targets -contains $text) {\n    [void]$toDelete.Add($p)\n  }\n}\n\n# Delete from the end backwards so earlier paragraph ranges stay valid.\nfor ($i = $toDelete.Count - 1; $i -ge 0; $i--) {\n  $toDelete[$i].Range.Delete()\n}\n"}
